# PortfolioData.xlsx refactor:
#   - Drop the "Prezzo Operazione" column's old date/live-price wiring
#     (OpenBB/OpenFIGI auto price lookups removed).
#   - Rename existing column H header to "Prezzo Operazione (EUR)".
#   - Add a new column I "Prezzo Corrente (EUR)" for manually entered
#     current prices.
#   - Clear out the auto-populated "Acquisto" rows (2 and 5) back to "-"
#     with blank operation price / date, since those were fed by the old
#     automated pipeline.
#   - Populate row 14 (EURIZON OB E HY) as a real manual "Acquisto" entry:
#     operation date, operation price (=+E14) and a manually typed current
#     price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----------------------------------------------------------
$ws.Range("H1").Value = "Prezzo Operazione (EUR)"
$ws.Range("I1").Value = "Prezzo Corrente (EUR)"

# ---- Row 2 (ACMB IL HEALTH C EUR) : revert to an untouched "-" row -------
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = $null

# ---- Row 5 (BPA CC E28) : revert to an untouched "-" row -----------------
$ws.Range("F5").Value = $null
$ws.Range("G5").Value = "-"
$ws.Range("H5").Value = $null
$ws.Range("I5").Value = $null

# ---- Row 14 (EURIZON OB E HY) : now a real manual "Acquisto" row --------
$ws.Range("F14").Value = "01/01/2024"
$ws.Range("F14").NumberFormat = "mm-dd-yy"
$ws.Range("G14").Value = "Acquisto"
$ws.Range("H14").Formula = "=+E14"
$ws.Range("I14").Value = 10.5

# ---- Column widths / new column I formatting -----------------------------
$ws.Columns("H").ColumnWidth = 21.36328125
$ws.Columns("I").ColumnWidth = 21.90625

# ---- Active cell back on I1, matching the author's saved selection ------
$ws.Range("I1").Select()
